$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row (F1, G1) ---
$ws.Range("F1").Value = "GI_Agilent"
$ws.Range("F1").Style = "Normal"
$ws.Range("G1").Value = "GI_Agilent_alternatif"
$ws.Range("G1").Style = "Normal"

# --- Data rows: F = numeric GI value, G = "NA" unless an alternative value exists ---
$fValues = @{
    2 = 9.14
    3 = 15
    4 = 4
    5 = 9.1
    6 = 20.6
    7 = 24.1
    8 = 1
    9 = 2
    10 = 3
    11 = 36
    12 = 39.2
    13 = 60
    14 = 9.14
    15 = 10.6
    16 = 10.6
    17 = 12.8
    18 = 13.5
    19 = 9.14
    20 = 7.2
    21 = 1
}

$gNumericValues = @{
    3 = 22.5
    18 = 8.16
}

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 6).Value = $fValues[$r]
    $ws.Cells.Item($r, 6).Style = "Normal"

    if ($gNumericValues.ContainsKey($r)) {
        $ws.Cells.Item($r, 7).Value = $gNumericValues[$r]
    } else {
        $ws.Cells.Item($r, 7).Value = "NA"
    }
    $ws.Cells.Item($r, 7).Style = "Normal"
}

# --- Restore default (non-custom) row heights on rows that previously had an explicit height ---
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(20).AutoFit()
$ws.Rows.Item(21).AutoFit()

# --- New note row below the table ---
$ws.Range("D31").Value = "S"

# --- Column G width (best-fit sized for the longer header text) ---
$ws.Columns.Item(7).ColumnWidth = 18.86

# --- Selection / active cell, matching the final authored state ---
$ws.Range("D31").Select()
